# update RFLP plate maps
$wb = $excel.ActiveWorkbook

# Rename first sheet "PocHistone RLFP 004" -> "DONE PocHistone RLFP 004"
$wsFirst = $wb.Worksheets.Item("PocHistone RLFP 004")
$wsFirst.Name = "DONE PocHistone RLFP 004"

# Move the active/selected tab from the first sheet to the last sheet
# ("PocHistone RLFP 019"), which also scrolls the tab strip so later
# sheets (starting around the 10th tab) are in view.
$wsLast = $wb.Worksheets.Item("PocHistone RLFP 019")
$wsLast.Activate()
